$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values (Fecha, Calidad, Volumen, Precio minimo, Precio maximo, Precio promedio ponderado, Precio $/Kg)
# for rows 2..18 as per the weekly refresh of the "Hortaliza - Ciboulette" dataset.
$rows = @(
    @{ Row = 2;  D = 45149; I = "Primera"; J = 80;  K = 2500; L = 2500; M = 2500; P = 833 },
    @{ Row = 3;  D = 45149; I = "Segunda"; J = 80;  K = 2000; L = 2000; M = 2000; P = 667 },
    @{ Row = 4;  D = 45146; I = "Primera"; J = 80;  K = 2500; L = 2500; M = 2500; P = 833 },
    @{ Row = 5;  D = 45146; I = "Segunda"; J = 80;  K = 2000; L = 2000; M = 2000; P = 667 },
    @{ Row = 6;  D = 44838; I = "Primera"; J = 200; K = 1200; L = 1300; M = 1250; P = 417 },
    @{ Row = 7;  D = 44838; I = "Segunda"; J = 150; K = 1000; L = 1000; M = 1000; P = 333 },
    @{ Row = 8;  D = 45148; I = "Primera"; J = 80;  K = 2500; L = 2500; M = 2500; P = 833 },
    @{ Row = 9;  D = 45148; I = "Segunda"; J = 60;  K = 2000; L = 2000; M = 2000; P = 667 },
    @{ Row = 10; D = 45133; I = "Primera"; J = 80;  K = 2500; L = 2500; M = 2500; P = 833 },
    @{ Row = 11; D = 44846; I = "Primera"; J = 200; K = 1200; L = 1300; M = 1250; P = 417 },
    @{ Row = 12; D = 44846; I = "Segunda"; J = 150; K = 1000; L = 1000; M = 1000; P = 333 },
    @{ Row = 13; D = 44832; I = "Primera"; J = 200; K = 1200; L = 1300; M = 1250; P = 417 },
    @{ Row = 14; D = 44832; I = "Segunda"; J = 150; K = 1000; L = 1000; M = 1000; P = 333 },
    @{ Row = 15; D = 45145; I = "Primera"; J = 60;  K = 2500; L = 2500; M = 2500; P = 833 },
    @{ Row = 16; D = 45145; I = "Segunda"; J = 80;  K = 2000; L = 2000; M = 2000; P = 667 },
    @{ Row = 17; D = 45134; I = "Primera"; J = 50;  K = 2500; L = 2500; M = 2500; P = 833 },
    @{ Row = 18; D = 45135; I = "Primera"; J = 70;  K = 2500; L = 2500; M = 2500; P = 833 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Cells.Item($n, 4).Value  = $r.D   # D: Fecha (serial date)
    $ws.Cells.Item($n, 9).Value  = $r.I   # I: Calidad
    $ws.Cells.Item($n, 10).Value = $r.J   # J: Volumen
    $ws.Cells.Item($n, 11).Value = $r.K   # K: Precio minimo
    $ws.Cells.Item($n, 12).Value = $r.L   # L: Precio maximo
    $ws.Cells.Item($n, 13).Value = $r.M   # M: Precio promedio ponderado
    $ws.Cells.Item($n, 16).Value = $r.P   # P: Precio $/Kg
}
